$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add two new rows to the table (this grows the table range / autofilter automatically).
$row3 = $lo.ListRows.Add()
$row4 = $lo.ListRows.Add()

# Row 3 data: 2025-05-11, 19:00 -> 21:00, offline
$ws.Range("A3").Value = 45788
$ws.Range("B3").Value = 0.79166666666666663
$ws.Range("C3").Value = 0.875
$ws.Range("D3").Formula = "=Table1[[#This Row],[To]]-Table1[[#This Row],[From]]"
$ws.Range("E3").Value = "offline"

# Row 4 data: 2025-05-15, 18:30 -> 22:00, offline
$ws.Range("A4").Value = 45792
$ws.Range("B4").Value = 0.77083333333333337
$ws.Range("C4").Value = 0.91666666666666663
$ws.Range("D4").Formula = "=Table1[[#This Row],[To]]-Table1[[#This Row],[From]]"
$ws.Range("E4").Value = "offline"

# Switch the Date column to a long-date display (e.g. "Thursday, May 08, 2025").
$ws.Range("A2").NumberFormat = "[$-F800]dddd, mmmm dd, yyyy"

# Copy the (now fully-formatted) first data row down onto the two new rows, so every
# column keeps the same number format / alignment as row 2.
$ws.Range("A2:E2").Copy()
$ws.Range("A3:E4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column A widens to fit the new, longer date text (closest achievable match to the
# ~20.33 "characters" width Excel would auto-fit to for the long-date format).
$ws.Columns.Item(1).ColumnWidth = 19.5

# Keep the selection where the user left off after entering the new data.
$ws.Range("E5").Select()
